$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(44488, 40, 16000, 17000, 16500, 1269)
    3  = @(44425, 60, 14000, 15000, 14500, 1115)
    5  = @(44327, 50, 24000, 25000, 24400, 1877)
    6  = @(44523, 40, 15000, 16000, 15500, 1192)
    7  = @(44308, 50, 26000, 27000, 26400, 2031)
    8  = @(44453, 50, 14000, 15000, 14600, 1123)
    9  = @(44377, 40, 14000, 15000, 14500, 1115)
    10 = @(44383, 50, 15000, 16000, 15400, 1185)
    11 = @(44435, 100, 13000, 14000, 13500, 1038)
    12 = @(44467, 100, 13000, 14000, 13500, 1038)
    13 = @(44510, 40, 15000, 16000, 15500, 1192)
    14 = @(44474, 40, 13000, 14000, 13500, 1038)
    15 = @(44350, 40, 23000, 25000, 24000, 1846)
    16 = @(44509, 100, 15000, 16000, 15500, 1192)
    17 = @(44433, 100, 13000, 14000, 13500, 1038)
    18 = @(44316, 50, 27000, 28000, 27400, 2108)
    19 = @(44159, 60, 30000, 32000, 31000, 2385)
    20 = @(44320, 50, 26000, 28000, 26800, 2062)
    21 = @(44313, 50, 25000, 26000, 25600, 1969)
    22 = @(44334, 50, 26000, 28000, 27200, 2092)
    23 = @(44362, 40, 15000, 16000, 15500, 1192)
    24 = @(44462, 60, 14000, 15000, 14500, 1115)
    25 = @(44355, 60, 18000, 20000, 19000, 1462)
    26 = @(44264, 40, 30000, 32000, 31000, 2385)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
    $ws.Cells.Item($row, 11).Value = $vals[2]
    $ws.Cells.Item($row, 12).Value = $vals[3]
    $ws.Cells.Item($row, 13).Value = $vals[4]
    $ws.Cells.Item($row, 16).Value = $vals[5]
}
